# Updates cryptos list figures (prices / 1h volume change) for Sheet1, rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "26.149.92"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").Value = "1.658.33"
$ws.Range("E3").Value = "  -0.62%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'217.26"
$ws.Range("E5").Value = "  -1.17%  "

$ws.Range("D6").Value = "'0.5217"
$ws.Range("E6").Value = "  -0.71%  "

$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.2635"
$ws.Range("E8").Value = "  -1.26%  "

$ws.Range("D9").Value = "'0.06261"
$ws.Range("E9").Value = "  -1.46%  "

$ws.Range("D10").Value = "'20.69"
$ws.Range("E10").Value = "  -4.45%  "

$ws.Range("D11").Value = "'0.07758"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.674.12"
$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.453"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").Value = "1.888.21"
$ws.Range("E14").Value = "  -0.54%  "

$ws.Range("D15").Value = "'0.5452"
$ws.Range("E15").Value = "  -1.11%  "

$ws.Range("D16").Value = "0.0₅8114"
$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("D17").Value = "'64.81"
$ws.Range("E17").Value = "  -1.06%  "

$ws.Range("D18").Value = "26.178.96"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("D20").Value = "'4.576"
$ws.Range("E20").Value = "  -3.26%  "

$ws.Range("D21").Value = "'191.53"
$ws.Range("E21").Value = "  -0.87%  "

$ws.Range("D22").Value = "'9.995"
$ws.Range("E22").Value = "  -2.57%  "

$ws.Range("D23").Value = "'5.987"
$ws.Range("E23").Value = "  -4.27%  "

$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").Value = "'138.39"
$ws.Range("E25").Value = "  -0.46%  "

$ws.Range("D26").Value = "'0.1229"
$ws.Range("E26").Value = "  -2.66%  "

$ws.Range("E27").Value = "  -1.60%  "

$ws.Range("D28").Value = "'16.18"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").Value = "'1.403"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("D30").Value = "'0.05929"
$ws.Range("E30").Value = "  -3.35%  "

$ws.Range("D31").Value = "'1.276"
$ws.Range("E31").Value = "  -1.11%  "

$ws.Range("D32").Value = "'3.524"
$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("D33").Value = "'3.250"
$ws.Range("E33").Value = "  -4.13%  "

$ws.Range("D34").Value = "'1.571"
$ws.Range("E34").Value = "  -6.12%  "

$ws.Range("D35").Value = "'0.9543"
$ws.Range("E35").Value = "  -4.65%  "

$ws.Range("D36").Value = "'2.418"

$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("D38").Value = "'0.5668"
$ws.Range("E38").Value = "  -5.96%  "

$ws.Range("E39").Value = "  -0.88%  "

$ws.Range("D40").Value = "'5.938"
$ws.Range("E40").Value = "  -1.55%  "

$ws.Range("D41").Value = "'0.8478"
$ws.Range("E41").Value = "  -0.93%  "

$ws.Range("D42").Value = "'1.002"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "'100.49"
$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("D44").Value = "1.002.77"
$ws.Range("E44").Value = "  -7.64%  "

$ws.Range("D45").Value = "1.802.69"
$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("D46").Value = "'56.42"
$ws.Range("E46").Value = "  -2.37%  "

$ws.Range("E47").Value = "  -3.78%  "

$ws.Range("E48").Value = "  -0.29%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.4345"
$ws.Range("E49").Value = "  +2.70%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'7.974"
$ws.Range("E50").Value = "  -2.34%  "

$ws.Range("D51").Value = "'0.05154"
$ws.Range("E51").Value = "  -0.91%  "
